$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date-serial value for every data row
# (rows 2 through 206). All of them are being bumped forward by one day
# (46074 -> 46075), so increment each cell's existing numeric value by 1.
for ($row = 2; $row -le 206; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
